# Update the "Förändrad" (Changed) date column (C) for all data rows
# from the old serial date value 45190 (2023-09-21) to the new value
# 45192 (2023-09-23), matching the committed change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 72; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
